$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11:B12").NumberFormat = "@"

$ws.Range("A11").Value = "Chris"
$ws.Range("B11").Value = "000987"
$ws.Range("A12").Value = "Jonny"
$ws.Range("B12").Value = "111000"

$ws.Range("A1").Select()
